$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2-3: account holder / card number
$ws.Range("C2").Value = "Hartmut"
# B3 holds a 16-digit card number that must stay text (not become a
# number / lose precision to scientific notation). Writing it via a
# formula that evaluates to a text string, then freezing the result
# with Copy + PasteSpecial(Values), keeps the cell's existing style
# while avoiding Excel's automatic numeric coercion of .Value.
$ws.Range("B3").Formula = "=""2570314725427075"""
$ws.Range("B3").Copy()
$ws.Range("B3").PasteSpecial(-4163)
$ws.Range("C3").Value = "Mohaupt"

# Row 5: opening balance date
$ws.Range("D5").Value = "KONTOSTAND AM 16.06.2024"

# Row 6: transaction 1
$ws.Range("B6").Value = "18.06."
$ws.Range("C6").Value = "19.06."
$ws.Range("D6").Value = "KARTENZAHLUNG ARAL TANKSTELLE"
$ws.Range("E6").Value = "60,97-"

# Row 7: transaction 2
$ws.Range("B7").Value = "20.06."
$ws.Range("C7").Value = "21.06."
$ws.Range("D7").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E7").Value = "24,94-"

# Row 8: transaction 3
$ws.Range("B8").Value = "24.06."
$ws.Range("C8").Value = "25.06."
$ws.Range("D8").Value = "RECHNUNG VODAFONE GMBH 1794266"
$ws.Range("E8").Value = "40,41-"

# Row 9: transaction 4
$ws.Range("B9").Value = "26.06."
$ws.Range("C9").Value = "27.06."
$ws.Range("D9").Value = "BEITRAG Allianz SE K-51295472"
$ws.Range("E9").Value = "56,98-"

# Row 10: transaction removed - clear cells entirely (becomes blank row,
# matching the blank-row style pattern already used by row 11). E10's
# alignment is set to match the target style (right/center/wrap) before
# clearing its value so the cell keeps a non-default style (style 12)
# instead of reverting to the default style when emptied.
$ws.Range("B10").Value = ""
$ws.Range("C10").Value = ""
$ws.Range("D10").Value = ""
$ws.Range("E10").HorizontalAlignment = -4152
$ws.Range("E10").VerticalAlignment = -4108
$ws.Range("E10").WrapText = $true
$ws.Range("E10").Value = ""

# Row 12: closing balance date and amount
$ws.Range("D12").Value = "KONTOSTAND AM 30.06.2024"
$ws.Range("E12").Value = "183,30-"

# Row 13: next statement date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 06.07.2024"
